{"js": "// Replace the date line and every \"a\u00f7b=\" prompt in the practice table with the\n// new values from the source commit. Every old string is unique in the document,\n// so a single forward search/replace pass (one hit each) is used per pair.\n// The pass runs top-to-bottom in document order: the cell that originally read\n// \"47\u00f76=\" is updated to \"61\u00f72=\" before a *different* cell (\"14\u00f72=\") is updated to\n// the text \"47\u00f76=\", so the freshly written \"47\u00f76=\" is never matched/touched again.\nconst replacements = [\n  [\"2025-07-05 Saturday\", \"2025-07-06 Sunday\"],\n  [\"63\u00f77=\", \"82\u00f76=\"],\n  [\"79\u00f77=\", \"53\u00f74=\"],\n  [\"39\u00f75=\", \"60\u00f77=\"],\n  [\"51\u00f74=\", \"76\u00f72=\"],\n  [\"47\u00f76=\", \"61\u00f72=\"],\n  [\"93\u00f78=\", \"97\u00f73=\"],\n  [\"48\u00f77=\", \"17\u00f77=\"],\n  [\"27\u00f78=\", \"89\u00f74=\"],\n  [\"18\u00f73=\", \"96\u00f76=\"],\n  [\"47\u00f73=\", \"51\u00f72=\"],\n  [\"16\u00f76=\", \"40\u00f74=\"],\n  [\"65\u00f78=\", \"90\u00f77=\"],\n  [\"31\u00f72=\", \"21\u00f72=\"],\n  [\"69\u00f73=\", \"67\u00f77=\"],\n  [\"49\u00f74=\", \"32\u00f72=\"],\n  [\"80\u00f73=\", \"54\u00f72=\"],\n  [\"14\u00f72=\", \"47\u00f76=\"],\n  [\"99\u00f74=\", \"70\u00f72=\"],\n  [\"26\u00f72=\", \"86\u00f77=\"],\n  [\"53\u00f73=\", \"25\u00f74=\"],\n  [\"26\u00f75=\", \"62\u00f79=\"],\n  [\"11\u00f78=\", \"30\u00f79=\"],\n  [\"68\u00f74=\", \"18\u00f72=\"],\n  [\"64\u00f78=\", \"85\u00f74=\"],\n  [\"21\u00f78=\", \"67\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"a\u00f7b=\" prompt in the practice table with the\n# new values from the source commit. Every old string is unique in the document,\n# so a single forward Find/Replace pass (one hit each, not \"replace all\") is used.\n# The pass must run top-to-bottom in document order: the cell that originally read\n# \"47\u00f76=\" is updated to \"61\u00f72=\" before a *different* cell (\"14\u00f72=\") is updated to\n# the text \"47\u00f76=\", so the freshly written \"47\u00f76=\" is never touched again.\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"2025-07-05 Saturday\",\n    \"63\u00f77=\",\n    \"79\u00f77=\",\n    \"39\u00f75=\",\n    \"51\u00f74=\",\n    \"47\u00f76=\",\n    \"93\u00f78=\",\n    \"48\u00f77=\",\n    \"27\u00f78=\",\n    \"18\u00f73=\",\n    \"47\u00f73=\",\n    \"16\u00f76=\",\n    \"65\u00f78=\",\n    \"31\u00f72=\",\n    \"69\u00f73=\",\n    \"49\u00f74=\",\n    \"80\u00f73=\",\n    \"14\u00f72=\",\n    \"99\u00f74=\",\n    \"26\u00f72=\",\n    \"53\u00f73=\",\n    \"26\u00f75=\",\n    \"11\u00f78=\",\n    \"68\u00f74=\",\n    \"64\u00f78=\",\n    \"21\u00f78=\"\n)\n\n$newValues = @(\n    \"2025-07-06 Sunday\",\n    \"82\u00f76=\",\n    \"53\u00f74=\",\n    \"60\u00f77=\",\n    \"76\u00f72=\",\n    \"61\u00f72=\",\n    \"97\u00f73=\",\n    \"17\u00f77=\",\n    \"89\u00f74=\",\n    \"96\u00f76=\",\n    \"51\u00f72=\",\n    \"40\u00f74=\",\n    \"90\u00f77=\",\n    \"21\u00f72=\",\n    \"67\u00f77=\",\n    \"32\u00f72=\",\n    \"54\u00f72=\",\n    \"47\u00f76=\",\n    \"70\u00f72=\",\n    \"86\u00f77=\",\n    \"25\u00f74=\",\n    \"62\u00f79=\",\n    \"30\u00f79=\",\n    \"18\u00f72=\",\n    \"85\u00f74=\",\n    \"67\u00f73=\"\n)\n\nfor ($i = 0; $i -lt $oldValues.Count; $i++) {\n    $rng = $d.Content\n    $rng.Find.Execute($oldValues[$i], $false, $false, $false, $false, $false, $true, 1, $false, $newValues[$i], 2, $false) | Out-Null\n}\n"}
